# Auto-generated Excel COM-interop edit script
# Applies the cryptos.xlsx price/volume/coin update described in the commit diff
# ("Updated cryptos list on Sun Oct 20 13:40:05 UTC 2024 with GitHub Actions").

function Set-TextValue($Range, $Text) {
    # Force the cell to be treated as TEXT so numeric-looking strings
    # (e.g. "597.18") are kept verbatim as strings instead of being
    # auto-coerced into (imprecise) floating point numbers, while still
    # leaving the cell with no explicit style afterwards (matches source).
    $Range.NumberFormat = "@"
    $Range.Value = $Text
    $Range.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$v = '68.416.49'
Set-TextValue $ws.Range("D2") $v
$ws.Range("E2").Value = '  +0.41%  '

$v = '2.652.38'
Set-TextValue $ws.Range("D3") $v
$ws.Range("E3").Value = '  +0.59%  '

$ws.Range("E4").Value = '  -0.07%  '

$v = '597.18'
Set-TextValue $ws.Range("D5") $v
$ws.Range("E5").Value = '  +0.09%  '

$v = '158.20'
Set-TextValue $ws.Range("D6") $v
$ws.Range("E6").Value = '  +2.22%  '

$ws.Range("E7").Value = '  -0.03%  '

$v = '0.540'
Set-TextValue $ws.Range("D8") $v
$ws.Range("E8").Value = '  -0.97%  '

$v = '2.650.15'
Set-TextValue $ws.Range("D9") $v
$ws.Range("E9").Value = '  +0.56%  '

$v = '0.139'
Set-TextValue $ws.Range("D10") $v
$ws.Range("E10").Value = '  -3.85%  '

$ws.Range("E11").Value = '  -0.82%  '

$v = '5.27'
Set-TextValue $ws.Range("D12") $v
$ws.Range("E12").Value = '  +0.58%  '

$v = '0.355'
Set-TextValue $ws.Range("D13") $v
$ws.Range("E13").Value = '  +1.40%  '

$v = '28.07'
Set-TextValue $ws.Range("D14") $v
$ws.Range("E14").Value = '  +0.78%  '

$v = '3.131.81'
Set-TextValue $ws.Range("D15") $v
$ws.Range("E15").Value = '  +0.42%  '

$v = '0.0000186'
Set-TextValue $ws.Range("D16") $v
$ws.Range("E16").Value = '  -3.25%  '

$v = '68.326.21'
Set-TextValue $ws.Range("D17") $v
$ws.Range("E17").Value = '  +0.37%  '

$v = '2.642.13'
Set-TextValue $ws.Range("D18") $v
$ws.Range("E18").Value = '  +0.50%  '

$v = '11.74'
Set-TextValue $ws.Range("D19") $v
$ws.Range("E19").Value = '  +3.34%  '

$v = '364.10'
Set-TextValue $ws.Range("D20") $v
$ws.Range("E20").Value = '  +0.41%  '

$v = '7.55'
Set-TextValue $ws.Range("D21") $v
$ws.Range("E21").Value = '  +1.70%  '

$ws.Range("E22").Value = '  +2.43%  '

$v = '4.82'
Set-TextValue $ws.Range("D23") $v
$ws.Range("E23").Value = '  +0.32%  '

$v = '2.08'
Set-TextValue $ws.Range("D24") $v
$ws.Range("E24").Value = '  +1.26%  '

$v = '75.12'
Set-TextValue $ws.Range("D25") $v
$ws.Range("E25").Value = '  +0.28%  '

$ws.Range("E28").Value = '  +0.49%  '

$ws.Range("E29").Value = '  -2.66%  '

$ws.Range("E30").Value = '  -0.09%  '

$v = '575.25'
Set-TextValue $ws.Range("D31") $v
$ws.Range("E31").Value = '  +2.62%  '

$v = '8.23'
Set-TextValue $ws.Range("D32") $v
$ws.Range("E32").Value = '  +3.11%  '

$v = '1.42'
Set-TextValue $ws.Range("D33") $v
$ws.Range("E33").Value = '  +1.58%  '

$ws.Range("E34").Value = '  +2.52%  '

$ws.Range("E35").Value = '  +5.14%  '

$v = '160.80'
Set-TextValue $ws.Range("D38") $v
$ws.Range("E38").Value = '  -0.25%  '

$v = '19.69'
Set-TextValue $ws.Range("D39") $v
$ws.Range("E39").Value = '  +2.09%  '

$v = '0.374'
Set-TextValue $ws.Range("D40") $v
$ws.Range("E40").Value = '  +0.55%  '

$v = '1.89'
Set-TextValue $ws.Range("D41") $v
$ws.Range("E41").Value = '  +0.90%  '

$v = '5.33'
Set-TextValue $ws.Range("D42") $v
$ws.Range("E42").Value = '  +0.69%  '

$ws.Range("E43").Value = '  -0.85%  '

$ws.Range("E44").Value = '  +0.06%  '

$v = ('0.0{0}0315' -f [char]0x2086)
Set-TextValue $ws.Range("D45") $v
$ws.Range("E45").Value = '  -7.45%  '

$v = '158.64'
Set-TextValue $ws.Range("D46") $v
$ws.Range("E46").Value = '  +0.35%  '

$ws.Range("E47").Value = '  +4.43%  '

$ws.Range("E48").Value = '  +3.17%  '

$v = '0.0783'
Set-TextValue $ws.Range("D51") $v
$ws.Range("E51").Value = '  -0.22%  '

# Rows 49 and 50 swapped rank position (ARBITRUM <-> InjectiveProtocol)
$ws.Range("B49").Value = 'InjectiveProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$v = '21.95'
Set-TextValue $ws.Range("D49") $v
$ws.Range("E49").Value = '  +0.26%  '

$ws.Range("B50").Value = 'ARBITRUM'
$ws.Range("C50").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$v = '0.593'
Set-TextValue $ws.Range("D50") $v
$ws.Range("E50").Value = '  +6.23%  '
